$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15: fill in actual time + difference formula ---
$ws.Range("C15").Value = 45
$ws.Range("D15").Formula = "=B15-C15"

# --- Row 16: Self Evaluation: SWOT (still "Heading 2" look, like rows 6-15) ---
$ws.Range("A16").Style = "Heading 2"
$ws.Range("A16").Value = "Self Evaluation: SWOT"
$ws.Range("B16").Value = 60
$ws.Range("C16").Value = 70
$ws.Range("D16").Formula = "=B16-C16"

# --- Row 17: WK 2: Anchor Points ---
$ws.Range("A17").Style = "Heading 2"
$ws.Range("A17").Value = "WK 2: Anchor Points"
$ws.Range("B17").Value = 30
$ws.Range("C17").Value = 45
$ws.Range("D17").Formula = "=B17-C17"

# --- Row 18: WK 2: Failure to Success ---
$ws.Range("A18").Style = "Heading 2"
$ws.Range("A18").Value = "WK 2: Failure to Success"
$ws.Range("B18").Value = 90
$ws.Range("C18").Value = 120
$ws.Range("D18").Formula = "=B18-C18"

# --- Row 19: Mission Statement ---
$ws.Range("A19").Style = "Heading 2"
$ws.Range("A19").Value = "Mission Statement"
$ws.Range("B19").Value = 30
$ws.Range("C19").Value = 50
$ws.Range("D19").Formula = "=B19-C19"

# --- Row 20: WK 2: Project & Portfolio ---
$ws.Range("A20").Style = "Heading 2"
$ws.Range("A20").Value = "WK 2: Project & Portfolio"
$ws.Range("B20").Value = 120
$ws.Range("C20").Value = 150
$ws.Range("D20").Formula = "=B20-C20"

# --- Convert D7:D12 into one shared-formula block (re-enter the same relative formula) ---
$ws.Range("D7:D12").Formula = "=B7-C7"

# --- Trailing blank row under the new table ---
$ws.Rows.Item(21).RowHeight = 17

# --- Selection moves to where the user ended up typing ---
$ws.Range("D21").Select()
